$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Promociones")

# Update the placeholder text in D4 from NombreListaPrecio to ListaPrecio
$ws.Range("D4").Value = "{{item.ListaPrecio}}"

# Restore the active cell selection to J1 as in the edited workbook
$ws.Range("J1").Select()
